{"js": "// The document is a title paragraph (\"2025-03-16 Sunday\") followed by a\n// 20x5 table of arithmetic expressions (\"80-27=\", \"40+17=\", ...). The\n// commit regenerates every expression (and bumps the date by one day),\n// while keeping the overall paragraph count (1 title + 100 cells = 101)\n// and table shape (20 rows x 5 columns) unchanged.\n//\n// `context.document.body.paragraphs` walks the document in reading\n// order \u2014 the title paragraph first, then each table cell's paragraph\n// row-by-row, column-by-column \u2014 so we can apply the new values\n// positionally without needing to separately touch the table object.\n\nconst newTexts = [\"2025-03-17 Monday\", \"86-71=\", \"38-31=\", \"60-17=\", \"22+28=\", \"74-40=\", \"61-21=\", \"98-84=\", \"18+62=\", \"84+12=\", \"59+17=\", \"28+11=\", \"78-62=\", \"73-67=\", \"87+2=\", \"91-67=\", \"97-89=\", \"46+30=\", \"58-51=\", \"85+8=\", \"3+75=\", \"66-29=\", \"53-39=\", \"13+85=\", \"34+38=\", \"5+75=\", \"17+25=\", \"39+3=\", \"77-23=\", \"33+48=\", \"98-46=\", \"96-35=\", \"56+29=\", \"71-64=\", \"42-24=\", \"14+61=\", \"96-69=\", \"49-0=\", \"54+20=\", \"18+28=\", \"39+35=\", \"73-70=\", \"39-10=\", \"72-7=\", \"75+10=\", \"37+5=\", \"48-8=\", \"93-18=\", \"37+8=\", \"91-14=\", \"96-95=\", \"92-89=\", \"71+21=\", \"34+55=\", \"98-45=\", \"93-52=\", \"41+24=\", \"89-67=\", \"91-90=\", \"11+59=\", \"18+56=\", \"31-11=\", \"10+62=\", \"6+75=\", \"69-53=\", \"42+38=\", \"80-20=\", \"1+51=\", \"72+15=\", \"95-9=\", \"16+63=\", \"3+73=\", \"51+48=\", \"21+3=\", \"66-27=\", \"69-28=\", \"60-9=\", \"60-9=\", \"99-14=\", \"34+8=\", \"71-71=\", \"38-7=\", \"52-30=\", \"58-14=\", \"12+1=\", \"49-21=\", \"57+22=\", \"1+2=\", \"36+38=\", \"7+42=\", \"84-12=\", \"30-23=\", \"98-74=\", \"52-35=\", \"6+23=\", \"84-13=\", \"18+50=\", \"42+34=\", \"5+60=\", \"50+2=\", \"61+32=\"];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newTexts.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + newTexts.length +\n    \" got \" + paragraphs.items.length\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const current = paragraphs.items[i];\n  current.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const current = paragraphs.items[i];\n  if (current.text !== newTexts[i]) {\n    // Replace-in-place keeps the existing run formatting (font, size, ...).\n    current.insertText(newTexts[i], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document is a title paragraph (\"2025-03-16 Sunday\") followed by a\n# 20x5 table of arithmetic expressions (\"80-27=\", \"40+17=\", ...). The\n# commit regenerates every expression (and bumps the date by one day)\n# while keeping the table shape (20 rows x 5 columns) unchanged.\n#\n# NOTE: $d.Paragraphs.Count reports 121 (not 101) here because the COM\n# paragraph-enumeration model inserts a virtual/empty paragraph at each\n# table row boundary; those extra entries are not backed by real content\n# (writing to them is a no-op) so we avoid Paragraphs for the table and\n# instead address cells directly through the Tables/Cell API, which maps\n# 1:1 onto the <w:tc> cells in the OOXML.\n\n$d = $word.ActiveDocument\n\n# 1) Update the title paragraph (the single non-table paragraph).\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.Text = \"2025-03-17 Monday\"\n\n# 2) Update every table cell, row by row, left to right.\n$newCells = @(\n    @(\"86-71=\", \"38-31=\", \"60-17=\", \"22+28=\", \"74-40=\"),\n    @(\"61-21=\", \"98-84=\", \"18+62=\", \"84+12=\", \"59+17=\"),\n    @(\"28+11=\", \"78-62=\", \"73-67=\", \"87+2=\", \"91-67=\"),\n    @(\"97-89=\", \"46+30=\", \"58-51=\", \"85+8=\", \"3+75=\"),\n    @(\"66-29=\", \"53-39=\", \"13+85=\", \"34+38=\", \"5+75=\"),\n    @(\"17+25=\", \"39+3=\", \"77-23=\", \"33+48=\", \"98-46=\"),\n    @(\"96-35=\", \"56+29=\", \"71-64=\", \"42-24=\", \"14+61=\"),\n    @(\"96-69=\", \"49-0=\", \"54+20=\", \"18+28=\", \"39+35=\"),\n    @(\"73-70=\", \"39-10=\", \"72-7=\", \"75+10=\", \"37+5=\"),\n    @(\"48-8=\", \"93-18=\", \"37+8=\", \"91-14=\", \"96-95=\"),\n    @(\"92-89=\", \"71+21=\", \"34+55=\", \"98-45=\", \"93-52=\"),\n    @(\"41+24=\", \"89-67=\", \"91-90=\", \"11+59=\", \"18+56=\"),\n    @(\"31-11=\", \"10+62=\", \"6+75=\", \"69-53=\", \"42+38=\"),\n    @(\"80-20=\", \"1+51=\", \"72+15=\", \"95-9=\", \"16+63=\"),\n    @(\"3+73=\", \"51+48=\", \"21+3=\", \"66-27=\", \"69-28=\"),\n    @(\"60-9=\", \"60-9=\", \"99-14=\", \"34+8=\", \"71-71=\"),\n    @(\"38-7=\", \"52-30=\", \"58-14=\", \"12+1=\", \"49-21=\"),\n    @(\"57+22=\", \"1+2=\", \"36+38=\", \"7+42=\", \"84-12=\"),\n    @(\"30-23=\", \"98-74=\", \"52-35=\", \"6+23=\", \"84-13=\"),\n    @(\"18+50=\", \"42+34=\", \"5+60=\", \"50+2=\", \"61+32=\")\n)\n\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newCells.Count; $r++) {\n    $rowValues = $newCells[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
